$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: mark the run that holds the inline picture as "no proof"
# (adds <w:rPr><w:noProof/></w:rPr> inside the run that contains the
# <w:drawing> element).
# -----------------------------------------------------------------
$pic = $d.InlineShapes.Item(1)
$pic.Range.NoProofing = 1

# -----------------------------------------------------------------
# Change 2: fill in the answers for Precision / Recall.
# The document ends with two empty "Normal (Web)" paragraphs. The
# very last paragraph gets the "Precision: ..." text (two runs, the
# first one carrying a lastRenderedPageBreak marker), and a brand
# new paragraph with the same style is appended after it holding the
# "Recall: ..." text (also two runs).
# -----------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $lastPara.Range.Start

$precisionXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Precision: </w:t></w:r><w:r><w:t>"Out of all the points predicted to be positive, how many of them were actually positive?"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r1 = $d.Range($insertPos, $insertPos)
[void]$r1.InsertXML($precisionXml)

$endPos = $d.Content.End

$recallXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="120" w:beforeAutospacing="0"/></w:pPr><w:r><w:t xml:space="preserve">Recall: </w:t></w:r><w:r><w:t>"Out of the points that are labeled positive, how many of them were correctly predicted as positive?"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r2 = $d.Range($endPos, $endPos)
[void]$r2.InsertXML($recallXml)
